$d = $word.ActiveDocument

# --- Update the date line at the top of the document ---
$dateParagraph = $d.Paragraphs.Item(1)
$dateParagraph.Range.Text = "2024-12-29 Sunday"

# --- Update each division problem in the table ---
# Cells are addressed by (row, column) and the text is assigned directly
# to the cell Range, which correctly scopes the edit to that single cell
# (important since some problem values, e.g. "15÷7=", repeat elsewhere in the table).
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "83÷6="
$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "49÷8="
$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "43÷5="
$cell = $tbl.Cell(1, 4)
$cell.Range.Text = "93÷6="
$cell = $tbl.Cell(1, 5)
$cell.Range.Text = "50÷6="
$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "34÷9="
$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "88÷4="
$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "86÷7="
$cell = $tbl.Cell(5, 4)
$cell.Range.Text = "76÷9="
$cell = $tbl.Cell(5, 5)
$cell.Range.Text = "29÷4="
$cell = $tbl.Cell(9, 1)
$cell.Range.Text = "90÷9="
$cell = $tbl.Cell(9, 2)
$cell.Range.Text = "97÷2="
$cell = $tbl.Cell(9, 3)
$cell.Range.Text = "84÷7="
$cell = $tbl.Cell(9, 4)
$cell.Range.Text = "53÷3="
$cell = $tbl.Cell(9, 5)
$cell.Range.Text = "78÷6="
$cell = $tbl.Cell(13, 1)
$cell.Range.Text = "80÷9="
$cell = $tbl.Cell(13, 2)
$cell.Range.Text = "37÷6="
$cell = $tbl.Cell(13, 3)
$cell.Range.Text = "13÷5="
$cell = $tbl.Cell(13, 4)
$cell.Range.Text = "82÷5="
$cell = $tbl.Cell(13, 5)
$cell.Range.Text = "82÷6="
$cell = $tbl.Cell(17, 1)
$cell.Range.Text = "74÷6="
$cell = $tbl.Cell(17, 2)
$cell.Range.Text = "82÷9="
$cell = $tbl.Cell(17, 3)
$cell.Range.Text = "81÷7="
$cell = $tbl.Cell(17, 4)
$cell.Range.Text = "91÷9="
$cell = $tbl.Cell(17, 5)
$cell.Range.Text = "71÷5="
